$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped from 45180
# (2023-09-11) to 45181 (2023-09-12) for every data row (rows 2 through 238).
$ws.Range("C2:C238").Value = 45181
